# Add a new worksheet "MegaMenuInfo" right after the existing sheet, containing
# megamenu link data, and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after the first (and currently only) sheet.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "MegaMenuInfo"

# Populate cell values in left-to-right, top-to-bottom order so that newly
# created shared-string entries land in the same order as the source workbook.
$newSheet.Range("A1").Value = "Path"
$newSheet.Range("B1").Value = "LinkName"
$newSheet.Range("C1").Value = "Language"
$newSheet.Range("A2").Value = "/"
$newSheet.Range("A3").Value = "/espanol/tipos"
$newSheet.Range("C2").Value = "english"
$newSheet.Range("C3").Value = "spanish"

# Header row re-uses the same bold/fill header style as sheet1's header row.
# Copy formatting (not values) from that cell over to the new header row.
$ws1.Range("A1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths for the two text columns.
$newSheet.Columns.Item(1).ColumnWidth = 20.6
$newSheet.Columns.Item(2).ColumnWidth = 23.3

# Sheet1's old selection (A14) is replaced by a simple header-row selection.
$ws1.Range("A1:C1").Select()

# The new sheet becomes the active / tab-selected sheet, with A4 selected
# (just past the data that was entered).
$newSheet.Activate()
$newSheet.Range("A4").Select()
